$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rows 91 and 92: swap all data (columns B..AC), keep column A (the
#        running index) fixed per row. ---
$row91 = $ws.Range("B91:AC91").Value2
$row92 = $ws.Range("B92:AC92").Value2

$ws.Range("B91:AC91").Value2 = $row92
$ws.Range("B92:AC92").Value2 = $row91

# --- 2) Append a new row 231 with a new match record. ---
# Copy the formatting (styles) from row 230 so the new row matches the
# existing look (bold/border/center on col A, date format on col E, etc.)
$ws.Range("A230:AC230").Copy()
$ws.Range("A231:AC231").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A231").Value = 229
$ws.Range("B231").Value = 7641728
$ws.Range("C231").Value = "Mexico Liga de Expansion"
$ws.Range("D231").Value = "Mexico Liga de Expansion"
$ws.Range("E231").Value = 45395.83680555555
$ws.Range("F231").Value = "Oaxaca"
$ws.Range("G231").Value = "Atletico Morelia"
$ws.Range("H231").Value = 3
$ws.Range("I231").Value = 3
$ws.Range("J231").Value = "D"
$ws.Range("K231").Value = 2.25
$ws.Range("L231").Value = 3.5
$ws.Range("M231").Value = 2.7
$ws.Range("N231").Value = 2.3
$ws.Range("O231").Value = 3.6
$ws.Range("P231").Value = 2.875
$ws.Range("Q231").Value = -0.25
$ws.Range("R231").Value = 2.025
$ws.Range("S231").Value = 1.775
$ws.Range("T231").Value = 2.75
$ws.Range("U231").Value = 1.9
$ws.Range("V231").Value = 1.9
$ws.Range("W231").Value = -1
$ws.Range("X231").Value = 2.6
$ws.Range("Y231").Value = -1
$ws.Range("Z231").Value = -0.5
$ws.Range("AA231").Value = 0.3875
$ws.Range("AB231").Value = 0.8999999999999999
$ws.Range("AC231").Value = -1
